# Register-help API added: new "Description" column inserted before the
# request-body column, endpoint paths updated to include ":userid", and
# three new rows describing the offered/offering-help endpoints appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: relocate the existing "request body (required)" column from
#     C to D by hand (this keeps column B/E formatting untouched and avoids
#     dragging the old, unused column-D width metadata along for the ride,
#     the way a blind Columns.Insert would). ---
$ws.Range("D1").Value = "request body (required)"

$ws.Range("D3").Value = "UserTypeId, email, password"

$ws.Range("D4").Value = "{userdetails: {}, addressdetails : {}, verificationdetails : {}"
$ws.Range("D4").WrapText = $true

$ws.Range("D5").Value = "UserTypeId, email, password"

$ws.Range("D6").Value = "{userdetails: {}, addressdetails : {}, verificationdetails : {}"
$ws.Range("D6").WrapText = $true

$ws.Range("D8").Value = "email,password"
$ws.Range("D8").WrapText = $true

# --- Step 2: clear out the old column C (content + formatting) so it can
#     be reused as the new "Description" column. ---
$ws.Range("C1:C8").Clear()

# --- Step 3: write the new Description column. ---
$ws.Range("C1").Value = "Description"
$ws.Range("C3").Value = "Create account for user"
$ws.Range("C4").Value = "update all details of user"
$ws.Range("C5").Value = "Create account for NGOs"
$ws.Range("C6").Value = "update all details of NGOs"
$ws.Range("C7").Value = "returns  list of  all active NGOs"
$ws.Range("C8").Value = "login for user or NGO"

# --- Step 4: update endpoint paths to take a :userid parameter. ---
$ws.Range("A4").Value = "/api/user/updatedetails/:userid"
$ws.Range("A6").Value = "/api/ngo/updatedetails/:userid"

# --- Step 5: append the new "offered help" / "offering help" API rows. ---
$ws.Range("A9").Value = "/api/offeredhelps"
$ws.Range("B9").Value = "get"
$ws.Range("C9").Value = "list all user offered helps "
$ws.Range("A9").Font.Size = 12

$ws.Range("A10").Value = "/api/user/offeredhelps/:userid"
$ws.Range("B10").Value = "get"
$ws.Range("C10").Value = "get all helps offered by the given user"
$ws.Range("A10").Font.Size = 12

$ws.Range("A11").Value = "/api/user/offeringhelp/:userid"
$ws.Range("B11").Value = "post"
$ws.Range("C11").Value = "register help"
$ws.Range("A11").Font.Size = 12

$ws.Rows("9").RowHeight = 15.75
$ws.Rows("10").RowHeight = 15.75
$ws.Rows("11").RowHeight = 15.75

# --- Step 6: widen the endpoint / description columns, and restore the
#     request-body column's original width now that it lives in D. ---
$ws.Columns("A").ColumnWidth = 31.307291666666668
$ws.Columns("C").ColumnWidth = 32.451822916666664
$ws.Columns("D").ColumnWidth = 37.736979166666664

# --- Step 7: move the active selection to D16, matching the saved view. ---
[void]$ws.Range("D16").Select()
